$d = $word.ActiveDocument

# 1. "position from 0 to 10000" -> "position from 0 to 72000"
$d.Content.Find.Execute("position from 0 to 10000 are the default values", $true, $false, $false, $false, $false, $true, 1, $false, "position from 0 to 72000 are the default values", 2) | Out-Null

# 2. "pics out" -> "picks out" (typo fix)
$d.Content.Find.Execute("and pics out what algorithm", $true, $false, $false, $false, $false, $true, 1, $false, "and picks out what algorithm", 2) | Out-Null

# 3. Table caption text update
$d.Content.Find.Execute("Top classes use bottom classes. Not actual class names, just descriptors", $true, $false, $false, $false, $false, $true, 1, $false, "Joint interface has a pointer to an instance of the two classes below it", 2) | Out-Null
